$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.162.44"
$ws.Range("E2").Value = "  +0.78%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.827.14"
$ws.Range("E3").Value = "  +0.51%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.012"
$ws.Range("E4").Value = "  +0.87%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.00"
$ws.Range("E5").Value = "  +1.02%  "
$ws.Range("E6").Value = "  +0.72%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4713"
$ws.Range("E7").Value = "  +0.81%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3649"
$ws.Range("E8").Value = "  -0.45%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07412"
$ws.Range("E9").Value = "  +0.93%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8819"
$ws.Range("E10").Value = "  +1.11%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.44"
$ws.Range("E11").Value = "  +0.97%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.905.37"
$ws.Range("E12").Value = "  +4.97%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07326"
$ws.Range("E13").Value = "  +3.00%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.392"
$ws.Range("E14").Value = "  -0.23%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "93.35"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.524"
$ws.Range("E16").Value = "  +0.26%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.009"
$ws.Range("E17").Value = "  +0.48%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008745"
$ws.Range("E18").Value = "  +0.42%  "
$ws.Range("E19").Value = "  +0.81%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "27.692.01"
$ws.Range("E20").Value = "  +2.68%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.66"
$ws.Range("E21").Value = "  +0.05%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.258"
$ws.Range("E22").Value = "  -0.65%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.59"
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.111.40"
$ws.Range("E24").Value = "  +3.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.884"
$ws.Range("E25").Value = "  -0.48%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.77"
$ws.Range("E26").Value = "  +0.60%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.56"
$ws.Range("E27").Value = "  +1.12%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.135"
$ws.Range("E28").Value = "  -0.77%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.190"
$ws.Range("E29").Value = "  -1.22%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "116.46"
$ws.Range("E30").Value = "  -0.53%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08947"
$ws.Range("E31").Value = "  +0.57%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.167"
$ws.Range("E32").Value = "  +0.50%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7429"
$ws.Range("E33").Value = "  -2.09%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.514"
$ws.Range("E34").Value = "  +0.40%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.942"
$ws.Range("E35").Value = "  +1.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.010"
$ws.Range("E36").Value = "  +0.80%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.089"
$ws.Range("E37").Value = "  -0.12%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05316"
$ws.Range("E38").Value = "  +0.48%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01947"
$ws.Range("E39").Value = "  +0.08%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.408"
$ws.Range("E40").Value = "  +1.42%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.946"
$ws.Range("E41").Value = "  -0.78%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.182"
$ws.Range("E42").Value = "  +0.24%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5264"
$ws.Range("E43").Value = "  -0.48%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1643"
$ws.Range("E44").Value = "  -0.53%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.424"
$ws.Range("E45").Value = "  -0.06%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4888"
$ws.Range("E46").Value = "  +0.42%  "
$ws.Range("E47").Value = "  -0.38%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.010"
$ws.Range("E48").Value = "  +0.81%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "104.75"
$ws.Range("E49").Value = "  +1.13%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.653"
$ws.Range("E50").Value = "  -0.71%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06301"
$ws.Range("E51").Value = "  +0.11%  "

Write-Host "Applied all cell updates"